$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mapping of row -> new F (dSF) value, per the diff
$values = @{
    2  = 1
    3  = -1
    4  = 2
    6  = 1
    7  = -4
    8  = -3
    9  = 1
    10 = -4
    11 = -2
    12 = 5
    13 = -1
    14 = 1
    15 = 1
    16 = 2
    18 = -2
    19 = 1
    21 = 5
    22 = 1
    23 = 1
    25 = 2
    26 = 5
    27 = -4
    28 = 6
    29 = -7
    30 = 9
    31 = -3
    33 = -2
    35 = 2
    36 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
